# Fix, and add base test.
# - Rename the "ExpectedMsg" header (column J) to "RetMsg" on both sheets.
# - Update the current cell selection on both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Module01")
$ws2 = $wb.Worksheets.Item("Module02")

# Rename header cell J1 from "ExpectedMsg" to "RetMsg" on both sheets.
$ws1.Range("J1").Value = "RetMsg"
$ws2.Range("J1").Value = "RetMsg"

# Update the selected cell on each sheet (cosmetic selection change).
$ws1.Activate()
$ws1.Range("H5").Select()

$ws2.Activate()
$ws2.Range("I6").Select()

# Restore the originally active sheet (Module01 is tabSelected="1").
$ws1.Activate()
